# Generate Report for Handoff
# Adds two new localization-status rows (for files
# d82698fd-a5d1-4630-9652-2d8c9cad5f7c and ff17d22a-23dd-4b3b-bd18-c9af7c3046a9)
# to the Overview sheet and to each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$guid1 = "d82698fd-a5d1-4630-9652-2d8c9cad5f7c"
$guid2 = "ff17d22a-23dd-4b3b-bd18-c9af7c3046a9"
$hash1 = "bbe7ebbbd88a64e89bec983c6eeb839898ea4405"
$hash2 = "2bae0e080b4e9f2494b324e73ba66f6981eb306a"

$commitA = "c89eeb42bf50e206a87017a099ca0d186a607644"
$commitB = "66cd11d8ea36d84bdd12b07dba03b56bbd7f6c98"

# -------------------------------------------------------------------
# Sheet "Overview": two new rows, one per new source file.
# Columns: A = File Name (hyperlink to the .md source),
#          B = zh-cn status, C = de-de status, D = Latest Handoff Date
# -------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitA + "/e2e/" + $guid1 + ".md", "", "", $guid1 + ".md")
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-27-12 10:27:59"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitB + "/e2e/" + $guid2 + ".md", "", "", $guid2 + ".md")
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-27-12 10:27:59"

# -------------------------------------------------------------------
# Locale sheets "zh-cn" and "de-de": two new rows each.
# Columns: A = Source File Name (hyperlink to .md)
#          B = File Extension (hyperlink, display ".md")
#          C = Status
#          D = Latest Handoff File (hyperlink to .xlf)
#          E = Latest Handoff Datetime
#          H = Latest Handback DateTime
#          I = Handoff Reason
# -------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Hash1File = $guid1 + "." + $hash1 + ".zh-cn.xlf"; Hash2File = $guid2 + "." + $hash2 + ".zh-cn.xlf"; HandoffDate = "2016-03-12 10:27:56" },
    @{ Name = "de-de"; Hash1File = $guid1 + "." + $hash1 + ".de-de.xlf"; Hash2File = $guid2 + "." + $hash2 + ".de-de.xlf"; HandoffDate = "2016-03-12 10:27:59" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)
    $localeName = $loc.Name

    # Row 4 -> guid1
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitA + "/e2e/" + $guid1 + ".md", "", "", $guid1 + ".md")
    $ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitA + "/e2e/" + $guid1 + ".md", "", "", ".md")
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000001/ol-handoff/OpenLocalizationTestOrg/oltest." + $localeName + "/ci/ht/" + $loc.Hash1File, "", "", $loc.Hash1File)
    $ws.Range("E4").Value = $loc.HandoffDate
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "Include"

    # Row 5 -> guid2
    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitB + "/e2e/" + $guid2 + ".md", "", "", $guid2 + ".md")
    $ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitB + "/e2e/" + $guid2 + ".md", "", "", ".md")
    $ws.Range("C5").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000002/ol-handoff/OpenLocalizationTestOrg/oltest." + $localeName + "/ci/ht/" + $loc.Hash2File, "", "", $loc.Hash2File)
    $ws.Range("E5").Value = $loc.HandoffDate
    $ws.Range("H5").Value = "0001-01-01 00:00:00"
    $ws.Range("I5").Value = "Include"
}

Write-Output "Report rows added for $guid1 and $guid2"
